$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: Morning Glass of Ether | Ether
$ws.Range("H15").Value = 249.02
$ws.Range("I15").Value = 249.02
$ws.Range("K15").Value = 747.0600000000001
$ws.Range("M15").Value = -578.0600000000001

# Row 29: Dripping with Venom | Weak Blinding Potion
$ws.Range("H29").Value = 3450
$ws.Range("J29").Value = 3000
$ws.Range("L29").Value = 9000
$ws.Range("N29").Value = -9562

# Row 33: Glazed and Confused | Clear Glass Lens
$ws.Range("H33").Value = 3368105.2
$ws.Range("I33").Value = 1207.7
$ws.Range("K33").Value = 1207.7
$ws.Range("M33").Value = -978.7

# Row 38: Just Give Him a Serum | Hi-Potion of Strength
$ws.Range("H38").Value = 2348.4546
$ws.Range("I38").Value = 610.8333
$ws.Range("J38").Value = 4433.6
$ws.Range("K38").Value = 1832.4999
$ws.Range("L38").Value = 13300.8
$ws.Range("M38").Value = -1460.4999
$ws.Range("N38").Value = -14044.8

# Row 74: Adhesive of Antipathy | Wing Glue
$ws.Range("H74").Value = 4171.4287
$ws.Range("I74").Value = 4171.4287
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 4171.4287
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -3235.4287
$ws.Range("N74").ClearContents()

# Row 77: It's Gonna Grow Back (L) | Wing Glue
$ws.Range("H77").Value = 4171.4287
$ws.Range("I77").Value = 4171.4287
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 20857.1435
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -16177.1435
$ws.Range("N77").ClearContents()

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2630.6233
$ws.Range("I138").Value = 1358.9286
$ws.Range("J138").Value = 3499.0977
$ws.Range("K138").Value = 4076.7858
$ws.Range("L138").Value = 10497.2931
$ws.Range("M138").Value = 1063.2142
$ws.Range("N138").Value = -20777.2931

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots | Bronze Ingot
$ws.Range("H2").Value = 1443.4166
$ws.Range("I2").Value = 1588.7142
$ws.Range("J2").Value = 1240
$ws.Range("K2").Value = 1588.7142
$ws.Range("L2").Value = 1240
$ws.Range("M2").Value = -1475.7142
$ws.Range("N2").Value = -1466

# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 1723.0344
$ws.Range("I74").Value = 1277.6666
$ws.Range("J74").Value = 2200.2144
$ws.Range("K74").Value = 1277.6666
$ws.Range("L74").Value = 2200.2144
$ws.Range("M74").Value = -403.6666
$ws.Range("N74").Value = -3948.2144

# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 1723.0344
$ws.Range("I77").Value = 1277.6666
$ws.Range("J77").Value = 2200.2144
$ws.Range("K77").Value = 6388.333000000001
$ws.Range("L77").Value = 11001.072
$ws.Range("M77").Value = -2020.333000000001
$ws.Range("N77").Value = -19737.072

# Row 110: Scheduled Maintenance | Deepgold Ingot
$ws.Range("H110").Value = 14099.667
$ws.Range("I110").Value = 18294
$ws.Range("J110").Value = 1516.6666
$ws.Range("K110").Value = 18294
$ws.Range("L110").Value = 1516.6666
$ws.Range("M110").Value = -16249
$ws.Range("N110").Value = -5606.6666

# Row 116: No Scope | Titanbronze Ingot
$ws.Range("H116").Value = 1443.4166
$ws.Range("I116").Value = 1588.7142
$ws.Range("J116").Value = 1240
$ws.Range("K116").Value = 1588.7142
$ws.Range("L116").Value = 1240
$ws.Range("M116").Value = 705.2858000000001
$ws.Range("N116").Value = -5828

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells | Bronze Ingot
$ws.Range("H3").Value = 1443.4166
$ws.Range("I3").Value = 1588.7142
$ws.Range("J3").Value = 1240
$ws.Range("K3").Value = 1588.7142
$ws.Range("L3").Value = 1240
$ws.Range("M3").Value = -1474.7142
$ws.Range("N3").Value = -1468

# Row 22: Riveting Run | Iron Rivets
$ws.Range("H22").Value = 523
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 634.5
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 634.5
$ws.Range("M22").Value = -127
$ws.Range("N22").Value = -980.5

# Row 133: Paring Is Caring | Mountain Chromite Hatchet
$ws.Range("H133").Value = 45724.223
$ws.Range("J133").Value = 45724.223
$ws.Range("L133").Value = 45724.223
$ws.Range("N133").Value = -55844.223

$ws = $wb.Worksheets.Item("CRP")
# Row 28: Militia on My Mind | Iron Lance
$ws.Range("H28").Value = 40000
$ws.Range("J28").Value = 40000
$ws.Range("L28").Value = 40000
$ws.Range("N28").Value = -40490

# Row 43: The Long Lance of the Law | Steel Halberd
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

# Row 86: Birch, Please | Birch Lumber
$ws.Range("H86").Value = 1937.1538
$ws.Range("I86").Value = 1885.5
$ws.Range("K86").Value = 1885.5
$ws.Range("M86").Value = -762.5

# Row 89: Built This City on Blocks and Soul (L) | Birch Lumber
$ws.Range("H89").Value = 1937.1538
$ws.Range("I89").Value = 1885.5
$ws.Range("K89").Value = 9427.5
$ws.Range("M89").Value = -3811.5

# Row 101: Everybody's Heard about the 'Berd | Doman Steel Halberd
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 68: Such a Butter Face | Fermented Butter
$ws.Range("H68").Value = 1887.6
$ws.Range("I68").Value = 458.66666
$ws.Range("J68").Value = 2500
$ws.Range("K68").Value = 1375.99998
$ws.Range("L68").Value = 7500
$ws.Range("M68").Value = -564.9999800000001
$ws.Range("N68").Value = -9122

# Row 71: No Margarine of Error (L) | Fermented Butter
$ws.Range("H71").Value = 1887.6
$ws.Range("I71").Value = 458.66666
$ws.Range("J71").Value = 2500
$ws.Range("K71").Value = 4127.99994
$ws.Range("L71").Value = 22500
$ws.Range("M71").Value = -71.9999399999997
$ws.Range("N71").Value = -30612

# Row 80: Saucy for a Suitor | Hollandaise Sauce
$ws.Range("H80").Value = 2000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 6000
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -7872

# Row 83: Saved by the Sauce (L) | Hollandaise Sauce
$ws.Range("H83").Value = 2000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 18000
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -27360

# Row 92: Oh No Udon | Gyr Abanian Flour
$ws.Range("H92").Value = 841.80646
$ws.Range("I92").Value = 768.44446
$ws.Range("J92").Value = 871.8182
$ws.Range("K92").Value = 2305.33338
$ws.Range("L92").Value = 2615.4546
$ws.Range("M92").Value = -1057.33338
$ws.Range("N92").Value = -5111.4546

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 1725291.5
$ws.Range("I131").Value = 4545972.5
$ws.Range("J131").Value = 1542.1111
$ws.Range("K131").Value = 13637917.5
$ws.Range("L131").Value = 4626.3333
$ws.Range("M131").Value = -13632877.5
$ws.Range("N131").Value = -14706.3333

# Row 134: Don't Knock It Till You've Tried It | Mezcal-marinated Swampmonk
$ws.Range("H134").Value = 12051.777
$ws.Range("I134").Value = 14741.5
$ws.Range("J134").Value = 9900
$ws.Range("K134").Value = 44224.5
$ws.Range("L134").Value = 29700
$ws.Range("M134").Value = -39154.5
$ws.Range("N134").Value = -39840

$ws = $wb.Worksheets.Item("GSM")
# Row 126: Gold Rush Order | Phrygian Gold Ingot
$ws.Range("H126").Value = 8940.444
$ws.Range("I126").Value = 11610.105
$ws.Range("J126").Value = 2600
$ws.Range("K126").Value = 34830.315
$ws.Range("L126").Value = 7800
$ws.Range("M126").Value = -32360.315
$ws.Range("N126").Value = -12740

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 3975.6445
$ws.Range("I132").Value = 3740.2334
$ws.Range("K132").Value = 11220.7002
$ws.Range("M132").Value = -8690.700199999999

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban | Leather
$ws.Range("H7").Value = 3095.5
$ws.Range("I7").Value = 2166.6667
$ws.Range("J7").Value = 4488.75
$ws.Range("K7").Value = 2166.6667
$ws.Range("L7").Value = 4488.75
$ws.Range("M7").Value = -2054.6667
$ws.Range("N7").Value = -4712.75

# Row 61: Spelling Me Softly | Raptor Leather
$ws.Range("H61").Value = 2587.2222
$ws.Range("I61").Value = 2341.4285
$ws.Range("J61").Value = 3447.5
$ws.Range("K61").Value = 2341.4285
$ws.Range("L61").Value = 3447.5
$ws.Range("M61").Value = -2139.4285
$ws.Range("N61").Value = -3851.5

# Row 113: Peace in Rest | Atrociraptor Leather
$ws.Range("H113").Value = 2587.2222
$ws.Range("I113").Value = 2341.4285
$ws.Range("J113").Value = 3447.5
$ws.Range("K113").Value = 2341.4285
$ws.Range("L113").Value = 3447.5
$ws.Range("M113").Value = -171.4285
$ws.Range("N113").Value = -7787.5

# Row 126: Battered Books | Saiga Leather
$ws.Range("H126").Value = 3095.5
$ws.Range("I126").Value = 2166.6667
$ws.Range("J126").Value = 4488.75
$ws.Range("K126").Value = 6500.000100000001
$ws.Range("L126").Value = 13466.25
$ws.Range("M126").Value = -4030.000100000001
$ws.Range("N126").Value = -18406.25

$ws = $wb.Worksheets.Item("WVR")
# Row 138: Halfgloves, Full Effort | Rroneek Serge Halfgloves of Healing
$ws.Range("H138").Value = 45002.4
$ws.Range("J138").Value = 41670.668
$ws.Range("L138").Value = 41670.668
$ws.Range("N138").Value = -51950.668
